$d = $word.ActiveDocument

# The hidden "_GoBack" bookmark currently sits at the end of paragraph 1
# (right after "Lab Assignment 5"). It needs to move to the start of the
# paragraph that begins with "An android application...". Remove it now
# so it doesn't "stick" to whatever content we insert at its old position,
# then re-add it at the correct spot once that paragraph's range is known.
$d.Bookmarks("_GoBack").Delete()

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Paragraph 1 ("Lab Assignment 5" heading) ---------------------------
# Collapse the 3rd, 4th and 5th <w:tab/> runs into a single run that keeps
# one tab and adds 17 literal spaces of text.
$p1 = $d.Paragraphs(1)
$p1TextEnd = $p1.Range.End - 1   # exclude the paragraph mark
$r1 = $d.Range($p1.Range.Start, $p1TextEnd)
$xml1 = "<w:p $ns>" +
          "<w:pPr>" +
            "<w:jc w:val=""both""/>" +
            "<w:rPr>" +
              "<w:rFonts w:asciiTheme=""majorHAnsi"" w:hAnsiTheme=""majorHAnsi""/>" +
              "<w:b/><w:sz w:val=""36""/><w:szCs w:val=""36""/>" +
            "</w:rPr>" +
          "</w:pPr>" +
          "<w:r><w:tab/></w:r>" +
          "<w:r><w:tab/></w:r>" +
          "<w:r><w:tab/><w:t xml:space=""preserve"">                 </w:t></w:r>" +
          "<w:r>" +
            "<w:rPr>" +
              "<w:rFonts w:asciiTheme=""majorHAnsi"" w:hAnsiTheme=""majorHAnsi""/>" +
              "<w:b/><w:sz w:val=""36""/><w:szCs w:val=""36""/>" +
            "</w:rPr>" +
            "<w:t>Lab Assignment 5</w:t>" +
          "</w:r>" +
        "</w:p>"
$r1.InsertXML($xml1)

# --- Paragraph 3 (the "An android application..." paragraph) ------------
# Drop the leading tab and the spell-check proofErr wrapper around "api",
# merging everything into one run, and put the _GoBack bookmark at the
# very start of the paragraph.
$p3 = $d.Paragraphs(3)
$p3TextEnd = $p3.Range.End - 1   # exclude the paragraph mark
$r3 = $d.Range($p3.Range.Start, $p3TextEnd)
$xml3 = "<w:p $ns>" +
          "<w:pPr><w:jc w:val=""both""/></w:pPr>" +
          "<w:bookmarkStart w:id=""0"" w:name=""_GoBack""/>" +
          "<w:bookmarkEnd w:id=""0""/>" +
          "<w:r><w:t>An android application which has a Registration page and map api which helps to find the location of the captured image.</w:t></w:r>" +
        "</w:p>"
$r3.InsertXML($xml3)
